# Fix test that started to fail after upgrading calc engine in ClosedXML;
# make the test culture-independent.
#
# The template cell B3 on sheet "Лист1" used a volatile, locale-dependent
# formula (TODAY()-7 formatted as "dd.MM.yyyy") whose cached/displayed text
# depends on the machine's current date and locale. Replace it with a
# static, culture-independent formula that concatenates a fixed date
# string, so the test's expected output no longer flaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("B3").Formula = '=CONCATENATE("Begin at ","19.01.2023")'
